# Update the AE detail table after the example accrual total changed
# from N=5 to N=3 (new start, version 0.0.10).
#
# This updates:
#   - the two column-header cells that mention "(N=5)" -> "(N=3)"
#   - every cell that displayed "20.00" (i.e. 1/5*100) -> "33.33" (1/3*100)
#   - every cell that displayed "40.00" (i.e. 2/5*100) -> "66.67" (2/3*100)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells: total accrual N=5 -> N=3 ---
$ws.Range("D9").Value = "the % of subjects that this comprises of the total accrual (N=3)"
$ws.Range("F9").Value = "% of the subjects that this comprises of the total accrual (N=3)"

# --- Cells showing "20.00" (1 of N) -> "33.33" (1 of 3) ---
# These values must stay stored as text (like the original "20.00" string),
# so the number format is forced to Text before assigning the value.
$cells20 = @("D10","F10","D11","F11","D12","D13","D15","D16","D17","D18","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","F35","D36","F36","D37","D38","D39","D41")
foreach ($addr in $cells20) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "33.33"
}

# --- Cells showing "40.00" (2 of N) -> "66.67" (2 of 3) ---
$cells40 = @("D14","D19","D21","D31","D40")
foreach ($addr in $cells40) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "66.67"
}
